$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.378.44'
$ws.Range('E2').Value = '  -2.08%  '
$ws.Range('D3').Value = '2.180.42'
$ws.Range('E3').Value = '  -3.07%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.09'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.54%  '
$ws.Range('E6').Value = '  -2.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.15'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.61%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.578'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.79'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.60%  '
$ws.Range('E11').Value = '  -3.28%  '
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.71'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.52%  '
$ws.Range('D14').Value = '2.509.35'
$ws.Range('E14').Value = '  -2.94%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.14'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.54%  '
$ws.Range('D16').Value = '2.183.72'
$ws.Range('E16').Value = '  -2.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.764'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.14%  '
$ws.Range('D18').Value = '42.279.07'
$ws.Range('E18').Value = '  -2.03%  '
$ws.Range('E19').Value = '  -4.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.82'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.94%  '
$ws.Range('E21').Value = '  -3.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '226.03'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.53%  '
$ws.Range('B23').Value = 'InternetComputer(DFINITY)'
$ws.Range('C23').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.26'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -13.63%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.09'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.61%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.38'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.39'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.62%  '
$ws.Range('E28').Value = '  -0.35%  '
$ws.Range('E29').Value = '  -4.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '37.18'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '171.53'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.99'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0819'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.16%  '
$ws.Range('E34').Value = '  -5.60%  '
$ws.Range('E35').Value = '  -2.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.106'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.74%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.15'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.96%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0332'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.05'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.72%  '
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.77'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -12.28%  '
$ws.Range('B41').Value = 'THORChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.14'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -8.42%  '
$ws.Range('E42').Value = '  -3.58%  '
$ws.Range('B43').Value = 'NEARProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.53'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +8.79%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '58.22'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.19'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.84%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0969'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.12%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.15'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.454'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.61%  '
$ws.Range('E49').Value = '  -3.65%  '
$ws.Range('E50').Value = '  -3.18%  '
$ws.Range('E51').Value = '  +0.00%  '
